$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 15
$ws.Range("H15").Value = 2187.5615
$ws.Range("I15").Value = 2187.5615
$ws.Range("K15").Value = 6562.684499999999
$ws.Range("M15").Value = -6393.684499999999

# ALC row 98
$ws.Range("H98").Value = 2468.65
$ws.Range("I98").Value = 2340.9473
$ws.Range("J98").Value = 4895
$ws.Range("K98").Value = 2340.9473
$ws.Range("L98").Value = 4895
$ws.Range("M98").Value = -842.9472999999998
$ws.Range("N98").Value = -7891

# ALC row 106
$ws.Range("H106").Value = 6311.778
$ws.Range("I106").Value = 7168.1665
$ws.Range("K106").Value = 7168.1665
$ws.Range("M106").Value = -6537.1665

# ALC row 112
$ws.Range("H112").Value = 3648.2632
$ws.Range("I112").Value = 1025
$ws.Range("J112").Value = 4347.8
$ws.Range("K112").Value = 3075
$ws.Range("L112").Value = 13043.4
$ws.Range("M112").Value = -1967
$ws.Range("N112").Value = -15259.4

# ALC row 122
$ws.Range("H122").Value = 2468.65
$ws.Range("I122").Value = 2340.9473
$ws.Range("J122").Value = 4895
$ws.Range("K122").Value = 7022.841899999999
$ws.Range("L122").Value = 14685
$ws.Range("M122").Value = -4572.841899999999
$ws.Range("N122").Value = -19585

# ALC row 132
$ws.Range("H132").Value = 16254.85
$ws.Range("I132").Value = 15228.286
$ws.Range("J132").Value = 18650.166
$ws.Range("K132").Value = 45684.858
$ws.Range("L132").Value = 55950.49800000001
$ws.Range("M132").Value = -43154.858
$ws.Range("N132").Value = -61010.49800000001

# ALC row 137
$ws.Range("H137").Value = 9890.866
$ws.Range("I137").Value = 4187.8
$ws.Range("J137").Value = 12742.4
$ws.Range("K137").Value = 12563.4
$ws.Range("L137").Value = 38227.2
$ws.Range("M137").Value = -10013.4
$ws.Range("N137").Value = -43327.2

# ALC row 138
$ws.Range("H138").Value = 2271.64
$ws.Range("I138").Value = 996.5454999999999
$ws.Range("J138").Value = 2631.282
$ws.Range("K138").Value = 2989.6365
$ws.Range("L138").Value = 7893.846
$ws.Range("M138").Value = 2150.3635
$ws.Range("N138").Value = -18173.846

# ALC row 141
$ws.Range("H141").Value = 731.65515
$ws.Range("I141").Value = 778.2593000000001
$ws.Range("J141").Value = 102.5
$ws.Range("K141").Value = 2334.7779
$ws.Range("L141").Value = 307.5
$ws.Range("M141").Value = 2845.2221
$ws.Range("N141").Value = -10667.5

$ws = $wb.Worksheets.Item("ARM")
# ARM row 2
$ws.Range("H2").Value = 955.5357
$ws.Range("I2").Value = 734.5238000000001
$ws.Range("J2").Value = 1618.5714
$ws.Range("K2").Value = 734.5238000000001
$ws.Range("L2").Value = 1618.5714
$ws.Range("M2").Value = -621.5238000000001
$ws.Range("N2").Value = -1844.5714

# ARM row 32
$ws.Range("H32").Value = 6758.75
$ws.Range("I32").Value = 1245.5513
$ws.Range("J32").Value = 26305.545
$ws.Range("K32").Value = 1245.5513
$ws.Range("L32").Value = 26305.545
$ws.Range("M32").Value = -958.5513000000001
$ws.Range("N32").Value = -26879.545

# ARM row 74
$ws.Range("H74").Value = 7546.028
$ws.Range("I74").Value = 1566.68
$ws.Range("J74").Value = 21135.455
$ws.Range("K74").Value = 1566.68
$ws.Range("L74").Value = 21135.455
$ws.Range("M74").Value = -692.6800000000001
$ws.Range("N74").Value = -22883.455

# ARM row 77
$ws.Range("H77").Value = 7546.028
$ws.Range("I77").Value = 1566.68
$ws.Range("J77").Value = 21135.455
$ws.Range("K77").Value = 7833.400000000001
$ws.Range("L77").Value = 105677.275
$ws.Range("M77").Value = -3465.400000000001
$ws.Range("N77").Value = -114413.275

# ARM row 116
$ws.Range("H116").Value = 955.5357
$ws.Range("I116").Value = 734.5238000000001
$ws.Range("J116").Value = 1618.5714
$ws.Range("K116").Value = 734.5238000000001
$ws.Range("L116").Value = 1618.5714
$ws.Range("M116").Value = 1559.4762
$ws.Range("N116").Value = -6206.5714

$ws = $wb.Worksheets.Item("BSM")
# BSM row 3
$ws.Range("H3").Value = 955.5357
$ws.Range("I3").Value = 734.5238000000001
$ws.Range("J3").Value = 1618.5714
$ws.Range("K3").Value = 734.5238000000001
$ws.Range("L3").Value = 1618.5714
$ws.Range("M3").Value = -620.5238000000001
$ws.Range("N3").Value = -1846.5714

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31
$ws.Range("H31").Value = 11329.667
$ws.Range("I31").Value = 4447.778
$ws.Range("J31").Value = 52621
$ws.Range("K31").Value = 4447.778
$ws.Range("L31").Value = 52621
$ws.Range("M31").Value = -4152.778
$ws.Range("N31").Value = -53211

# CRP row 34
$ws.Range("H34").Value = 11329.667
$ws.Range("I34").Value = 4447.778
$ws.Range("J34").Value = 52621
$ws.Range("K34").Value = 4447.778
$ws.Range("L34").Value = 52621
$ws.Range("M34").Value = -4245.778
$ws.Range("N34").Value = -53025

# CRP row 99
$ws.Range("H99").Value = 5595.9644
$ws.Range("I99").Value = 1576.5555
$ws.Range("J99").Value = 7499.8945
$ws.Range("K99").Value = 1576.5555
$ws.Range("L99").Value = 7499.8945
$ws.Range("M99").Value = -78.55549999999994
$ws.Range("N99").Value = -10495.8945

# CRP row 122
$ws.Range("H122").Value = 1979.125
$ws.Range("I122").Value = 1762.6666
$ws.Range("J122").Value = 2628.5
$ws.Range("K122").Value = 5287.9998
$ws.Range("L122").Value = 7885.5
$ws.Range("M122").Value = -2837.9998
$ws.Range("N122").Value = -12785.5

# CRP row 126
$ws.Range("H126").Value = 5595.9644
$ws.Range("I126").Value = 1576.5555
$ws.Range("J126").Value = 7499.8945
$ws.Range("K126").Value = 4729.666499999999
$ws.Range("L126").Value = 22499.6835
$ws.Range("M126").Value = -2259.666499999999
$ws.Range("N126").Value = -27439.6835

# CRP row 134
$ws.Range("H134").Value = 27784410
$ws.Range("I134").Value = 2019.8334
$ws.Range("J134").Value = 55566800
$ws.Range("K134").Value = 6059.5002
$ws.Range("L134").Value = 166700400
$ws.Range("M134").Value = -3524.5002
$ws.Range("N134").Value = -166705470

$ws = $wb.Worksheets.Item("CUL")
# CUL row 131
$ws.Range("H131").Value = 1453.94
$ws.Range("J131").Value = 1499.4255
$ws.Range("L131").Value = 4498.2765
$ws.Range("N131").Value = -14578.2765

# CUL row 132
$ws.Range("H132").Value = 1732.4445
$ws.Range("I132").Value = 1732.4445
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 15592.0005
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -13062.0005
$ws.Range("N132").ClearContents()

# CUL row 137
$ws.Range("H137").Value = 2592
$ws.Range("I137").Value = 2592
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 7776
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -2676
$ws.Range("N137").ClearContents()

# CUL row 140
$ws.Range("H140").Value = 1246.6666
$ws.Range("I140").Value = 1246.6666
$ws.Range("K140").Value = 3739.9998
$ws.Range("M140").Value = 1440.0002

$ws = $wb.Worksheets.Item("GSM")
# GSM row 122
$ws.Range("H122").Value = 5156.9473
$ws.Range("I122").Value = 5373.4443
$ws.Range("J122").Value = 1260
$ws.Range("K122").Value = 16120.3329
$ws.Range("L122").Value = 3780
$ws.Range("M122").Value = -13670.3329
$ws.Range("N122").Value = -8680

# GSM row 126
$ws.Range("H126").Value = 10490.429
$ws.Range("I126").Value = 15081.5
$ws.Range("J126").Value = 7665.154
$ws.Range("K126").Value = 45244.5
$ws.Range("L126").Value = 22995.462
$ws.Range("M126").Value = -42774.5
$ws.Range("N126").Value = -27935.462

$ws = $wb.Worksheets.Item("LTW")
# LTW row 93
$ws.Range("H93").Value = 16780.727
$ws.Range("I93").Value = 26498.5
$ws.Range("J93").Value = 11227.714
$ws.Range("K93").Value = 26498.5
$ws.Range("L93").Value = 11227.714
$ws.Range("M93").Value = -25250.5
$ws.Range("N93").Value = -13723.714

# LTW row 132
$ws.Range("H132").Value = 876401.3
$ws.Range("I132").Value = 1666.6875
$ws.Range("J132").Value = 2875794.8
$ws.Range("K132").Value = 5000.0625
$ws.Range("L132").Value = 8627384.399999999
$ws.Range("M132").Value = -2470.0625
$ws.Range("N132").Value = -8632444.399999999

$ws = $wb.Worksheets.Item("WVR")
# WVR row 96
$ws.Range("H96").Value = 2487.5
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 2487.5
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 2487.5
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -5233.5

# WVR row 122
$ws.Range("H122").Value = 2258.8823
$ws.Range("I122").Value = 960.28
$ws.Range("J122").Value = 5866.1113
$ws.Range("K122").Value = 2880.84
$ws.Range("L122").Value = 17598.3339
$ws.Range("M122").Value = -430.8400000000001
$ws.Range("N122").Value = -22498.3339

# WVR row 132
$ws.Range("H132").Value = 6047.3623
$ws.Range("I132").Value = 3432.0645
$ws.Range("J132").Value = 9050.111000000001
$ws.Range("K132").Value = 10296.1935
$ws.Range("L132").Value = 27150.333
$ws.Range("M132").Value = -7766.193499999999
$ws.Range("N132").Value = -32210.333

# WVR row 136
$ws.Range("H136").Value = 4539.5073
$ws.Range("I136").Value = 770.7234
$ws.Range("J136").Value = 13396.15
$ws.Range("K136").Value = 2312.1702
$ws.Range("L136").Value = 40188.45
$ws.Range("M136").Value = 237.8298
$ws.Range("N136").Value = -45288.45
